$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.658.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.760.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4493"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3735"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.65"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07794"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.208"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.381"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.760.61"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.14"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +12.66%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06265"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5311"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.694.63"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.338"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.83"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.349"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.960.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.22"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.215"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.782"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09285"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.695"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.79"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02339"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2187"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6506"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06139"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.040"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6009"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.83%  "
